# Applies the data edits described by the commit diff to the EnemyDB_Sheet workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Update stat values (rows 17-19, columns M/N/O)
$ws.Range("M17").Value = 10
$ws.Range("M18").Value = 10
$ws.Range("O18").Value = 5
$ws.Range("O19").Value = 6

# Move the active selection to M19, matching the saved view state in the diff
$ws.Range("M19").Select()
